# Scheduled market-data refresh: updates cached price/profit columns
# (H: currentAveragePrice, I: currentAveragePriceNQ, J: currentAveragePriceHQ,
#  K: LevePriceNQ, L: LevePriceHQ, M: LeveProfitNQ, N: LeveProfitHQ)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3184.5
$ws.Range("I76").Value = 2925.1428
$ws.Range("K76").Value = 2925.1428
$ws.Range("M76").Value = -2610.1428
$ws.Range("H79").Value = 3184.5
$ws.Range("I79").Value = 2925.1428
$ws.Range("K79").Value = 2925.1428
$ws.Range("M79").Value = -1833.1428
$ws.Range("H80").Value = 608.4167
$ws.Range("J80").Value = 700.1111
$ws.Range("L80").Value = 2100.3333
$ws.Range("N80").Value = -4096.3333
$ws.Range("H83").Value = 608.4167
$ws.Range("J83").Value = 700.1111
$ws.Range("L83").Value = 6300.9999
$ws.Range("N83").Value = -16284.9999
$ws.Range("H100").Value = 1567.6666
$ws.Range("I100").Value = 1787
$ws.Range("J100").Value = 800
$ws.Range("K100").Value = 1787
$ws.Range("L100").Value = 800
$ws.Range("M100").Value = -1246
$ws.Range("N100").Value = -1882
$ws.Range("H107").Value = 2232
$ws.Range("I107").Value = 2150.875
$ws.Range("J107").Value = 2313.125
$ws.Range("K107").Value = 2150.875
$ws.Range("L107").Value = 2313.125
$ws.Range("M107").Value = -230.875
$ws.Range("N107").Value = -6153.125
$ws.Range("H137").Value = 2188
$ws.Range("I137").Value = 1315.3334
$ws.Range("J137").Value = 2936
$ws.Range("K137").Value = 3946.0002
$ws.Range("L137").Value = 8808
$ws.Range("M137").Value = -1396.0002
$ws.Range("N137").Value = -13908
$ws.Range("H138").Value = 1954.78
$ws.Range("I138").Value = 772.6429000000001
$ws.Range("J138").Value = 2147.221
$ws.Range("K138").Value = 2317.9287
$ws.Range("L138").Value = 6441.663
$ws.Range("M138").Value = 2822.0713
$ws.Range("N138").Value = -16721.663

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3775.617
$ws.Range("I32").Value = 3699.4443
$ws.Range("K32").Value = 3699.4443
$ws.Range("M32").Value = -3412.4443
$ws.Range("H61").Value = 865.0454999999999
$ws.Range("I61").Value = 647.5263
$ws.Range("J61").Value = 2242.6667
$ws.Range("K61").Value = 647.5263
$ws.Range("L61").Value = 2242.6667
$ws.Range("M61").Value = -435.5263
$ws.Range("N61").Value = -2666.6667
$ws.Range("H74").Value = 1493.3
$ws.Range("I74").Value = 1325.8889
$ws.Range("J74").Value = 3000
$ws.Range("K74").Value = 1325.8889
$ws.Range("L74").Value = 3000
$ws.Range("M74").Value = -451.8888999999999
$ws.Range("N74").Value = -4748
$ws.Range("H77").Value = 1493.3
$ws.Range("I77").Value = 1325.8889
$ws.Range("J77").Value = 3000
$ws.Range("K77").Value = 6629.4445
$ws.Range("L77").Value = 15000
$ws.Range("M77").Value = -2261.4445
$ws.Range("N77").Value = -23736
$ws.Range("H97").Value = 494
$ws.Range("I97").Value = 449
$ws.Range("J97").Value = 606.5
$ws.Range("K97").Value = 449
$ws.Range("L97").Value = 606.5
$ws.Range("M97").Value = 47
$ws.Range("N97").Value = -1598.5
$ws.Range("H110").Value = 1458.5883
$ws.Range("I110").Value = 1118.8667
$ws.Range("K110").Value = 1118.8667
$ws.Range("M110").Value = 926.1333
$ws.Range("H132").Value = 2065.8572
$ws.Range("I132").Value = 1772.0952
$ws.Range("J132").Value = 3828.4285
$ws.Range("K132").Value = 5316.2856
$ws.Range("L132").Value = 11485.2855
$ws.Range("M132").Value = -2786.2856
$ws.Range("N132").Value = -16545.2855
$ws.Range("H136").Value = 865.0454999999999
$ws.Range("I136").Value = 647.5263
$ws.Range("J136").Value = 2242.6667
$ws.Range("K136").Value = 1942.5789
$ws.Range("L136").Value = 6728.000100000001
$ws.Range("M136").Value = 607.4211
$ws.Range("N136").Value = -11828.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 17858000
$ws.Range("I94").Value = 25000680
$ws.Range("J94").Value = 1299.75
$ws.Range("K94").Value = 25000680
$ws.Range("L94").Value = 1299.75
$ws.Range("M94").Value = -25000229
$ws.Range("N94").Value = -2201.75
$ws.Range("H134").Value = 5146.069
$ws.Range("I134").Value = 1731.7916
$ws.Range("K134").Value = 5195.3748
$ws.Range("M134").Value = -2660.3748

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2025780.8
$ws.Range("I99").Value = 3290942
$ws.Range("J99").Value = 1522.8
$ws.Range("K99").Value = 3290942
$ws.Range("L99").Value = 1522.8
$ws.Range("M99").Value = -3289444
$ws.Range("N99").Value = -4518.8
$ws.Range("H122").Value = 706.5
$ws.Range("I122").Value = 506
$ws.Range("K122").Value = 1518
$ws.Range("M122").Value = 932
$ws.Range("H126").Value = 2025780.8
$ws.Range("I126").Value = 3290942
$ws.Range("J126").Value = 1522.8
$ws.Range("K126").Value = 9872826
$ws.Range("L126").Value = 4568.4
$ws.Range("M126").Value = -9870356
$ws.Range("N126").Value = -9508.4
$ws.Range("H132").Value = 6224.074
$ws.Range("I132").Value = 7061.9
$ws.Range("K132").Value = 21185.7
$ws.Range("M132").Value = -18655.7
$ws.Range("H140").Value = 40000
$ws.Range("J140").Value = 40000
$ws.Range("L140").Value = 40000
$ws.Range("N140").Value = -50360
$ws.Range("H141").Value = 35243.332
$ws.Range("J141").Value = 35243.332
$ws.Range("L141").Value = 35243.332
$ws.Range("N141").Value = -45603.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 2000
$ws.Range("I25").Value = 2000
$ws.Range("J25").Value = 2000
$ws.Range("K25").Value = 6000
$ws.Range("L25").Value = 6000
$ws.Range("M25").Value = -5831
$ws.Range("N25").Value = -6338
$ws.Range("H30").Value = 2000
$ws.Range("I30").Value = 2000
$ws.Range("J30").Value = 2000
$ws.Range("K30").Value = 6000
$ws.Range("L30").Value = 6000
$ws.Range("M30").Value = -5898
$ws.Range("N30").Value = -6204
$ws.Range("H68").Value = 1989.9546
$ws.Range("I68").Value = 780.9091
$ws.Range("J68").Value = 3199
$ws.Range("K68").Value = 2342.7273
$ws.Range("L68").Value = 9597
$ws.Range("M68").Value = -1531.7273
$ws.Range("N68").Value = -11219
$ws.Range("H71").Value = 1989.9546
$ws.Range("I71").Value = 780.9091
$ws.Range("J71").Value = 3199
$ws.Range("K71").Value = 7028.1819
$ws.Range("L71").Value = 28791
$ws.Range("M71").Value = -2972.1819
$ws.Range("N71").Value = -36903
$ws.Range("H87").Value = 366.66666
$ws.Range("I87").Value = 366.66666
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 1099.99998
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = 148.0000199999999
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 366.66666
$ws.Range("I90").Value = 366.66666
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 3299.99994
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = 2940.00006
$ws.Range("N90").ClearContents()
$ws.Range("H99").Value = 1868.2307
$ws.Range("I99").Value = 468.75
$ws.Range("K99").Value = 1406.25
$ws.Range("M99").Value = 839.75
$ws.Range("H137").Value = 11290.167
$ws.Range("J137").Value = 15741.5
$ws.Range("L137").Value = 47224.5
$ws.Range("N137").Value = -57424.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 6400
$ws.Range("J29").Value = 6400
$ws.Range("L29").Value = 6400
$ws.Range("N29").Value = -6980
$ws.Range("H113").Value = 1268.75
$ws.Range("I113").Value = 1280
$ws.Range("J113").Value = 1235
$ws.Range("K113").Value = 1280
$ws.Range("L113").Value = 1235
$ws.Range("M113").Value = 890
$ws.Range("N113").Value = -5575

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 600
$ws.Range("I23").Value = 600
$ws.Range("K23").Value = 600
$ws.Range("M23").Value = -370
$ws.Range("H82").Value = 2286.6667
$ws.Range("I82").Value = 2207.5
$ws.Range("K82").Value = 2207.5
$ws.Range("M82").Value = -1846.5
$ws.Range("H85").Value = 2286.6667
$ws.Range("I85").Value = 2207.5
$ws.Range("K85").Value = 2207.5
$ws.Range("M85").Value = -959.5
$ws.Range("H122").Value = 11337048
$ws.Range("I122").Value = 20241518
$ws.Range("J122").Value = 4085.6365
$ws.Range("K122").Value = 60724554
$ws.Range("L122").Value = 12256.9095
$ws.Range("M122").Value = -60722104
$ws.Range("N122").Value = -17156.9095
$ws.Range("H132").Value = 50010.24
$ws.Range("I132").Value = 2109.7273
$ws.Range("J132").Value = 102700.8
$ws.Range("K132").Value = 6329.1819
$ws.Range("L132").Value = 308102.4
$ws.Range("M132").Value = -3799.1819
$ws.Range("N132").Value = -313162.4
$ws.Range("H136").Value = 6467.316
$ws.Range("I136").Value = 8019.2856
$ws.Range("J136").Value = 2121.8
$ws.Range("K136").Value = 24057.8568
$ws.Range("L136").Value = 6365.400000000001
$ws.Range("M136").Value = -21507.8568
$ws.Range("N136").Value = -11465.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 16875
$ws.Range("J92").Value = 16875
$ws.Range("L92").Value = 16875
$ws.Range("N92").Value = -21867
$ws.Range("H126").Value = 42736010
$ws.Range("I126").Value = 52910652
$ws.Range("J126").Value = 2500.2
$ws.Range("K126").Value = 158731956
$ws.Range("L126").Value = 7500.599999999999
$ws.Range("M126").Value = -158729486
$ws.Range("N126").Value = -12440.6
$ws.Range("H132").Value = 3867.7144
$ws.Range("I132").Value = 6021
$ws.Range("J132").Value = 1714.4286
$ws.Range("K132").Value = 18063
$ws.Range("L132").Value = 5143.2858
$ws.Range("M132").Value = -15533
$ws.Range("N132").Value = -10203.2858
$ws.Range("H141").Value = 43963
$ws.Range("J141").Value = 43963
$ws.Range("L141").Value = 43963
$ws.Range("N141").Value = -54323

